# Update Leave Card 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header date field (F3) now populated with a date
$ws.Range("F3").Value = 44743

# PERIOD column (A) shifted from "1st of month" to "end of month" for each
# bi-monthly/monthly leave-card row, plus newly-filled EARNED (C) entries
# for the rows that now have 1.25 credited.
$ws.Range("A11").Value = 44773
$ws.Range("A12").Value = 44804
$ws.Range("A13").Value = 44834
$ws.Range("A14").Value = 44865
$ws.Range("A15").Value = 44895
$ws.Range("A16").Value = 44926

$ws.Range("A18").Value = 44957
$ws.Range("A19").Value = 44985

$ws.Range("A20").Value = 45016
$ws.Range("C20").Value = 1.25

$ws.Range("A21").Value = 45046
$ws.Range("C21").Value = 1.25

$ws.Range("A22").Value = 45077
$ws.Range("C22").Value = 1.25

$ws.Range("A23").Value = 45107
$ws.Range("C23").Value = 1.25

$ws.Range("A24").Value = 45138
$ws.Range("C24").Value = 1.25

$ws.Range("A25").Value = 45169
$ws.Range("C25").Value = 1.25

$ws.Range("A26").Value = 45199
$ws.Range("C26").Value = 1.25

$ws.Range("A27").Value = 45230
$ws.Range("A28").Value = 45260
$ws.Range("A29").Value = 45291

# Update the cursor/selection so the saved view matches (bottom pane active
# cell moved from F5 to B15).
$ws.Range("B15").Select()
